$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "Report"
$ws.Name = "Report"

# Update membership counts (column B) for each society row
$ws.Range("B2").Value = 600
$ws.Range("B3").Value = 500
$ws.Range("B4").Value = 200
$ws.Range("B5").Value = 176
$ws.Range("B6").Value = 400
